# Update the "Volume / Number" line (A8) - issue number 12 -> 13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$volRange = $ws.Range("A8")
$volRange.Characters(21, 2).Text = "13"

# Update the "Report Covering the Week ... Through ..." line (C9)
# Week start: 3/20/2023 -> 3/27/2023 ; Week end: 3/26/2023 -> 4/2/2023
$weekRange = $ws.Range("C9")
$weekRange.Characters(27, 9).Text = "3/27/2023"
$weekRange.Characters(47, 9).Text = "4/2/2023"

# Row 14
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 6
$ws.Range("E14").Value = 16.666666666666
$ws.Range("F14").Value = 32
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = 6.666666666666
$ws.Range("I14").Value = 92
$ws.Range("J14").Value = 103
$ws.Range("K14").Value = -10.679611650485
$ws.Range("L14").Value = -11.538461538461
$ws.Range("M14").Value = -23.333333333333
$ws.Range("N14").Value = -81.300813008130

# Row 15
$ws.Range("C15").Value = 30
$ws.Range("D15").Value = 32
$ws.Range("E15").Value = -6.25
$ws.Range("F15").Value = 118
$ws.Range("G15").Value = 113
$ws.Range("H15").Value = 4.424778761061
$ws.Range("I15").Value = 375
$ws.Range("J15").Value = 410
$ws.Range("K15").Value = -8.536585365853
$ws.Range("L15").Value = 6.534090909090
$ws.Range("M15").Value = 25
$ws.Range("N15").Value = -50.396825396825

# Row 16
$ws.Range("C16").Value = 291
$ws.Range("D16").Value = 293
$ws.Range("E16").Value = -0.682593856655
$ws.Range("F16").Value = 1146
$ws.Range("G16").Value = 1160
$ws.Range("H16").Value = -1.206896551724
$ws.Range("I16").Value = 3843
$ws.Range("J16").Value = 3912
$ws.Range("K16").Value = -1.763803680981
$ws.Range("L16").Value = 44.256756756756
$ws.Range("M16").Value = -13.348365276212
$ws.Range("N16").Value = -82.027779076836

# Row 17
$ws.Range("C17").Value = 479
$ws.Range("D17").Value = 459
$ws.Range("E17").Value = 4.357298474945
$ws.Range("F17").Value = 1929
$ws.Range("G17").Value = 1871
$ws.Range("H17").Value = 3.099946552645
$ws.Range("I17").Value = 6220
$ws.Range("J17").Value = 5708
$ws.Range("K17").Value = 8.969866853538
$ws.Range("L17").Value = 32.115548003398
$ws.Range("M17").Value = 62.063574778530
$ws.Range("N17").Value = -30.495027377360

# Row 18
$ws.Range("C18").Value = 275
$ws.Range("D18").Value = 306
$ws.Range("E18").Value = -10.130718954248
$ws.Range("F18").Value = 1073
$ws.Range("G18").Value = 1236
$ws.Range("H18").Value = -13.187702265372
$ws.Range("I18").Value = 3627
$ws.Range("J18").Value = 3874
$ws.Range("K18").Value = -6.375838926174
$ws.Range("L18").Value = 23.032564450474
$ws.Range("M18").Value = -18.070928393946
$ws.Range("N18").Value = -85.354922070580

# Row 19
$ws.Range("C19").Value = 922
$ws.Range("D19").Value = 943
$ws.Range("E19").Value = -2.226935312831
$ws.Range("F19").Value = 3671
$ws.Range("G19").Value = 3685
$ws.Range("H19").Value = -0.379918588873
$ws.Range("I19").Value = 11927
$ws.Range("J19").Value = 12254
$ws.Range("K19").Value = -2.668516402807
$ws.Range("L19").Value = 53.461142563046
$ws.Range("M19").Value = 39.236516460424
$ws.Range("N19").Value = -38.374496228169

# Row 20
$ws.Range("C20").Value = 259
$ws.Range("D20").Value = 230
$ws.Range("E20").Value = 12.608695652173
$ws.Range("F20").Value = 1046
$ws.Range("G20").Value = 926
$ws.Range("H20").Value = 12.958963282937
$ws.Range("I20").Value = 3564
$ws.Range("J20").Value = 3331
$ws.Range("K20").Value = 6.994896427499
$ws.Range("L20").Value = 93.485342019544
$ws.Range("M20").Value = 50.952986022871
$ws.Range("N20").Value = -87.421916357861

# Row 21
$ws.Range("C21").Value = 2263
$ws.Range("D21").Value = 2269
$ws.Range("E21").Value = -0.264433671220
$ws.Range("F21").Value = 9015
$ws.Range("G21").Value = 9021
$ws.Range("H21").Value = -0.066511473229
$ws.Range("I21").Value = 29648
$ws.Range("J21").Value = 29592
$ws.Range("K21").Value = 0.189240335225
$ws.Range("L21").Value = 45.404610102991
$ws.Range("M21").Value = 23.291886721836
$ws.Range("N21").Value = -71.501898399577

# Row 22
$ws.Range("C22").Value = 54
$ws.Range("D22").Value = 42
$ws.Range("E22").Value = 28.571428571428
$ws.Range("F22").Value = 198
$ws.Range("G22").Value = 159
$ws.Range("H22").Value = 24.528301886792
$ws.Range("I22").Value = 540
$ws.Range("J22").Value = 587
$ws.Range("K22").Value = -8.006814310051
$ws.Range("L22").Value = 56.069364161849
$ws.Range("M22").Value = 5.058365758754

# Row 23
$ws.Range("C23").Value = 136
$ws.Range("D23").Value = 111
$ws.Range("E23").Value = 22.522522522522
$ws.Range("F23").Value = 461
$ws.Range("G23").Value = 420
$ws.Range("H23").Value = 9.761904761904
$ws.Range("I23").Value = 1503
$ws.Range("J23").Value = 1385
$ws.Range("K23").Value = 8.519855595667
$ws.Range("L23").Value = 25.041597337770
$ws.Range("M23").Value = 64.622124863088

# Row 24
$ws.Range("C24").Value = 2040
$ws.Range("D24").Value = 2147
$ws.Range("E24").Value = -4.983698183511
$ws.Range("F24").Value = 7992
$ws.Range("G24").Value = 8508
$ws.Range("H24").Value = -6.064880112834
$ws.Range("I24").Value = 26508
$ws.Range("J24").Value = 26077
$ws.Range("K24").Value = 1.652797484373
$ws.Range("L24").Value = 40.484392389633
$ws.Range("M24").Value = 45.736434108527

# Row 25
$ws.Range("C25").Value = 836
$ws.Range("D25").Value = 779
$ws.Range("E25").Value = 7.317073170731
$ws.Range("F25").Value = 3276
$ws.Range("G25").Value = 3172
$ws.Range("H25").Value = 3.278688524590
$ws.Range("I25").Value = 10146
$ws.Range("J25").Value = 9576
$ws.Range("K25").Value = 5.952380952380
$ws.Range("L25").Value = 37.071061875168
$ws.Range("M25").Value = -3.564299971485

# Row 26
$ws.Range("C26").Value = 46
$ws.Range("D26").Value = 45
$ws.Range("E26").Value = 2.222222222222
$ws.Range("F26").Value = 191
$ws.Range("G26").Value = 191
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 602
$ws.Range("J26").Value = 657
$ws.Range("K26").Value = -8.371385083713
$ws.Range("L26").Value = 3.793103448275

# Row 27
$ws.Range("C27").Value = 105
$ws.Range("D27").Value = 102
$ws.Range("E27").Value = 2.941176470588
$ws.Range("F27").Value = 416
$ws.Range("G27").Value = 415
$ws.Range("H27").Value = 0.240963855421
$ws.Range("I27").Value = 1234
$ws.Range("J27").Value = 1159
$ws.Range("K27").Value = 6.471095772217
$ws.Range("L27").Value = 24.020100502512

# Row 28
$ws.Range("C28").Value = 22
$ws.Range("D28").Value = 19
$ws.Range("E28").Value = 15.789473684210
$ws.Range("F28").Value = 97
$ws.Range("G28").Value = 121
$ws.Range("H28").Value = -19.834710743801
$ws.Range("I28").Value = 276
$ws.Range("J28").Value = 329
$ws.Range("K28").Value = -16.109422492401
$ws.Range("L28").Value = -3.496503496503
$ws.Range("M28").Value = -18.100890207715
$ws.Range("N28").Value = -79.956427015250

# Row 29
$ws.Range("C29").Value = 18
$ws.Range("D29").Value = 17
$ws.Range("E29").Value = 5.882352941176
$ws.Range("F29").Value = 81
$ws.Range("G29").Value = 106
$ws.Range("H29").Value = -23.584905660377
$ws.Range("I29").Value = 229
$ws.Range("J29").Value = 293
$ws.Range("K29").Value = -21.843003412969
$ws.Range("L29").Value = -10.894941634241
$ws.Range("M29").Value = -18.214285714285
$ws.Range("N29").Value = -81.825396825396

# Row 30
$ws.Range("C30").Value = 9
$ws.Range("D30").Value = 8
$ws.Range("E30").Value = 12.5
$ws.Range("F30").Value = 44
$ws.Range("G30").Value = 50
$ws.Range("H30").Value = -12
$ws.Range("I30").Value = 110
$ws.Range("J30").Value = 192
$ws.Range("K30").Value = -42.708333333333
$ws.Range("L30").Value = 23.595505617977

